$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: configure $p (a Word.Paragraph) as a ListParagraph bullet at the
# given (1-based) list level, reusing the document's existing list (numId 1),
# apply 1.5 "auto" line spacing, and set its text to a single run.
# ---------------------------------------------------------------------------
function Set-ListPara($p, $level, $text) {
    $p.Style = "List Paragraph"
    $p.Range.ListFormat.ListLevelNumber = $level
    $p.LineSpacingRule = 5
    $p.LineSpacing = 18
    $p.Range.Text = $text
}

# Same as Set-ListPara but without touching line spacing (keeps default).
function Set-ListParaNoSpacing($p, $level, $text) {
    $p.Style = "List Paragraph"
    $p.Range.ListFormat.ListLevelNumber = $level
    $p.Range.Text = $text
}

# ---------------------------------------------------------------------------
# Helper: append a brand-new paragraph at the end of the document and return
# it (inherits the pPr of whatever was previously last, which we then
# override as needed).
# ---------------------------------------------------------------------------
function Add-Para($doc) {
    $last = $doc.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    return $doc.Paragraphs.Last
}

# ---------------------------------------------------------------------------
# 1) The trailing blank paragraph already in the document becomes the new
#    "stackoverflow link" bullet (ilvl 1, numId 1).
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Last
Set-ListParaNoSpacing $p 2 "https://stackoverflow.com/questions/33229869/get-json-data-from-url-using-android"

# ---------------------------------------------------------------------------
# 2) Page break paragraph. Duplicate the styling/structure of an existing
#    lone page-break paragraph elsewhere in the document (copy/paste keeps
#    it a clean "<w:r><w:br w:type=\"page\"/></w:r>" with no stray pPr).
# ---------------------------------------------------------------------------
$p = Add-Para $d
$p.Style = "Normal"
$srcBreak = $d.Paragraphs(80).Range
$srcBreak.Copy()
$dst = $p.Range
$dst.Collapse(0)
$dst.Paste()

# ---------------------------------------------------------------------------
# 3) "Arduino computer software" Heading 1
# ---------------------------------------------------------------------------
$p = Add-Para $d
$p.Style = "Heading 1"
$p.LineSpacingRule = 5
$p.LineSpacing = 18
$p.Range.Text = "Arduino computer software"

# ---------------------------------------------------------------------------
# 4) "Packet structure" Heading 2
# ---------------------------------------------------------------------------
$p = Add-Para $d
$p.Style = "Heading 2"
$p.LineSpacingRule = 5
$p.LineSpacing = 18
$p.Range.Text = "Packet structure"

# ---------------------------------------------------------------------------
# 5) Bulleted structure description
# ---------------------------------------------------------------------------
$p = Add-Para $d
Set-ListPara $p 1 "Structure for packets from talkative devices to master:"

$p = Add-Para $d
Set-ListPara $p 2 "Sender ID (single alphanumerical character should suffice)"

$p = Add-Para $d
$p.Style = "List Paragraph"
$p.Range.ListFormat.ListLevelNumber = 2
$p.LineSpacingRule = 5
$p.LineSpacing = 18
$r = $p.Range
$r.InsertAfter("Packet urgency (")
$r.Collapse(0)
$r.InsertAfter("possible states: ")
$r.Collapse(0)
$r.InsertAfter("standard and error)")

$p = Add-Para $d
Set-ListPara $p 3 "This is used in case of power supply failover for example"

$p = Add-Para $d
Set-ListPara $p 3 "Forces master to listen and forward broadcast on network"

$p = Add-Para $d
Set-ListPara $p 3 "Standard packets update non-critical information"

$p = Add-Para $d
Set-ListPara $p 2 "Packet data"

$p = Add-Para $d
Set-ListPara $p 1 "Structure for packets from master to devices"

$p = Add-Para $d
Set-ListPara $p 2 "Target ID (alphanumerical character describing what computer is addressed)"

$p = Add-Para $d
$p.Style = "List Paragraph"
$p.Range.ListFormat.ListLevelNumber = 2
$p.LineSpacingRule = 5
$p.LineSpacing = 18
$r = $p.Range
$r.InsertAfter("Request type (")
$r.Collapse(0)
$r.InsertAfter("possible states:")
$r.Collapse(0)
$r.InsertAfter(" execute an action ")
$r.Collapse(0)
$r.InsertAfter("or")
$r.Collapse(0)
$r.InsertAfter(" request data)")

$p = Add-Para $d
Set-ListPara $p 2 "Request data (only necessary if asking a device to perform an action)"

$p = Add-Para $d
Set-ListPara $p 3 "Could be used for requesting oddly specific data that may be omitted from regular broadcasts"

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
